$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2135593220338983
$ws.Range("C2").Value = 0.5254237288135594
$ws.Range("J2").Value = 0.03389830508474576
$ws.Range("P2").Value = 0.1593220338983051
$ws.Range("S2").Value = 0.06779661016949153
$ws.Range("B3").Value = 0.02325581395348837
$ws.Range("C3").Value = 0.04651162790697674
$ws.Range("J3").Value = 0.08139534883720931
$ws.Range("P3").Value = 0.6686046511627907
$ws.Range("S3").Value = 0.1802325581395349
$ws.Range("J4").Value = 0.1578947368421053
$ws.Range("P4").Value = 0.5789473684210527
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.0759493670886076
$ws.Range("D6").Value = 0.008438818565400843
$ws.Range("F6").Value = 0.0970464135021097
$ws.Range("J6").Value = 0.1940928270042194
$ws.Range("O6").Value = 0.02109704641350211
$ws.Range("Q6").Value = 0.1772151898734177
$ws.Range("R6").Value = 0.06329113924050633
$ws.Range("S6").Value = 0.3628691983122363
$ws.Range("B7").Value = 0.1040462427745665
$ws.Range("D7").Value = 0.01734104046242774
$ws.Range("E7").Value = 0.005780346820809248
$ws.Range("F7").Value = 0.05780346820809248
$ws.Range("J7").Value = 0.1040462427745665
$ws.Range("Q7").Value = 0.2832369942196532
$ws.Range("R7").Value = 0.08670520231213873
$ws.Range("S7").Value = 0.3410404624277457
$ws.Range("B8").Value = 0.0774487471526196
$ws.Range("D8").Value = 0.01138952164009112
$ws.Range("F8").Value = 0.05694760820045558
$ws.Range("J8").Value = 0.1662870159453303
$ws.Range("O8").Value = 0.009111617312072893
$ws.Range("Q8").Value = 0.1776765375854214
$ws.Range("R8").Value = 0.1116173120728929
$ws.Range("S8").Value = 0.3895216400911162
$ws.Range("B9").Value = 0.07058823529411765
$ws.Range("D9").Value = 0.01764705882352941
$ws.Range("F9").Value = 0.05294117647058823
$ws.Range("J9").Value = 0.1764705882352941
$ws.Range("O9").Value = 0.01176470588235294
$ws.Range("Q9").Value = 0.1588235294117647
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.4117647058823529
$ws.Range("B10").Value = 0.0975609756097561
$ws.Range("D10").Value = 0.01672473867595819
$ws.Range("F10").Value = 0.07317073170731707
$ws.Range("J10").Value = 0.1554006968641115
$ws.Range("O10").Value = 0.01602787456445993
$ws.Range("Q10").Value = 0.2285714285714286
$ws.Range("R10").Value = 0.06689895470383275
$ws.Range("S10").Value = 0.3456445993031359
$ws.Range("G11").Value = 0.13671875
$ws.Range("J11").Value = 0.10546875
$ws.Range("K11").Value = 0.19140625
$ws.Range("L11").Value = 0.55859375
$ws.Range("S11").Value = 0.0078125
$ws.Range("G12").Value = 0.7960526315789473
$ws.Range("J12").Value = 0.125
$ws.Range("K12").Value = 0.006578947368421052
$ws.Range("L12").Value = 0.04605263157894737
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.282051282051282
$ws.Range("S13").Value = 0.05128205128205128
$ws.Range("F15").Value = 0.02403846153846154
$ws.Range("H15").Value = 0.1490384615384615
$ws.Range("I15").Value = 0.0673076923076923
$ws.Range("J15").Value = 0.4086538461538461
$ws.Range("K15").Value = 0.08653846153846154
$ws.Range("M15").Value = 0.004807692307692308
$ws.Range("O15").Value = 0.09134615384615384
$ws.Range("S15").Value = 0.1682692307692308
$ws.Range("F16").Value = 0.0115606936416185
$ws.Range("H16").Value = 0.2023121387283237
$ws.Range("I16").Value = 0.06358381502890173
$ws.Range("J16").Value = 0.4855491329479769
$ws.Range("K16").Value = 0.08092485549132948
$ws.Range("M16").Value = 0.03468208092485549
$ws.Range("O16").Value = 0.02890173410404624
$ws.Range("S16").Value = 0.09248554913294797
$ws.Range("F17").Value = 0.01346153846153846
$ws.Range("H17").Value = 0.1769230769230769
$ws.Range("I17").Value = 0.07884615384615384
$ws.Range("J17").Value = 0.4384615384615385
$ws.Range("K17").Value = 0.1
$ws.Range("M17").Value = 0.02115384615384616
$ws.Range("N17").Value = 0.003846153846153846
$ws.Range("O17").Value = 0.0673076923076923
$ws.Range("S17").Value = 0.1
$ws.Range("F18").Value = 0.02512562814070352
$ws.Range("H18").Value = 0.1658291457286432
$ws.Range("I18").Value = 0.06532663316582915
$ws.Range("J18").Value = 0.4673366834170855
$ws.Range("K18").Value = 0.09547738693467336
$ws.Range("M18").Value = 0.02010050251256281
$ws.Range("O18").Value = 0.06532663316582915
$ws.Range("S18").Value = 0.09547738693467336
$ws.Range("F19").Value = 0.01354784081287045
$ws.Range("H19").Value = 0.214225232853514
$ws.Range("I19").Value = 0.07705334462320068
$ws.Range("J19").Value = 0.4140558848433531
$ws.Range("K19").Value = 0.09144792548687553
$ws.Range("M19").Value = 0.01439458086367485
$ws.Range("O19").Value = 0.06689246401354784
$ws.Range("S19").Value = 0.1083827265029636
